$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '22.395.54'
$ws.Range("E2").Value = '  +0.06%  '

$ws.Range("D3").Value = '1.570.54'
$ws.Range("E3").Value = '  -0.09%  '

$ws.Range("E4").Value = '  +0.07%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '1.003'
$ws.Range("E5").Value = '  +0.16%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '291.61'
$ws.Range("E6").Value = '  +0.33%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3759'
$ws.Range("E7").Value = '  +2.48%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '49.75'
$ws.Range("E8").Value = '  +0.53%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.3411'
$ws.Range("E9").Value = '  +0.24%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.07617'
$ws.Range("E10").Value = '  -0.11%  '

$ws.Range("E11").Value = '  -2.07%  '

$ws.Range("E12").Value = '  -0.08%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '21.16'
$ws.Range("E13").Value = '  -0.93%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.008'
$ws.Range("E14").Value = '  -0.79%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.954'
$ws.Range("E15").Value = '  +0.38%  '

$ws.Range("D16").Value = '1.578.66'
$ws.Range("E16").Value = '  +0.57%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.00001133'
$ws.Range("E17").Value = '  -0.33%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '90.17'
$ws.Range("E18").Value = '  +0.84%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06749'
$ws.Range("E19").Value = '  +0.04%  '

$ws.Range("E20").Value = '  +0.04%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '16.74'
$ws.Range("E21").Value = '  +1.20%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.184'
$ws.Range("E22").Value = '  -0.97%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '11.98'
$ws.Range("E23").Value = '  -0.03%  '

$ws.Range("D24").Value = '22.390.62'
$ws.Range("E24").Value = '  +0.01%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.390'
$ws.Range("E25").Value = '  +0.36%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.684'
$ws.Range("E26").Value = '  -8.62%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '20.10'
$ws.Range("E27").Value = '  +0.48%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '147.28'
$ws.Range("E28").Value = '  +0.93%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '5.041'
$ws.Range("E29").Value = '  +1.41%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '126.66'
$ws.Range("E30").Value = '  +0.75%  '

$ws.Range("D31").Value = '1.750.13'
$ws.Range("E31").Value = '  +0.29%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.017'
$ws.Range("E32").Value = '  +0.65%  '

$ws.Range("B33").Value = 'ImmutableX'
$ws.Range("C33").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.001'
$ws.Range("E33").Value = '  -4.08%  '

$ws.Range("B34").Value = 'Filecoin'
$ws.Range("C34").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '6.099'
$ws.Range("E34").Value = '  -2.87%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '10.14'
$ws.Range("E35").Value = '  -1.21%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.08470'
$ws.Range("E36").Value = '  -0.51%  '

$ws.Range("B37").Value = 'TrustWalletToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.396'
$ws.Range("E37").Value = '  +11.47%  '

$ws.Range("B38").Value = 'VeChain'
$ws.Range("C38").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.02536'
$ws.Range("E38").Value = '  +0.30%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.2308'
$ws.Range("E39").Value = '  -0.80%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.06491'
$ws.Range("E40").Value = '  -0.81%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '5.425'
$ws.Range("E41").Value = '  -2.25%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '11.40'
$ws.Range("E42").Value = '  -2.81%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.6335'
$ws.Range("E43").Value = '  -0.52%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.002'
$ws.Range("E44").Value = '  +0.12%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '14.07'
$ws.Range("E45").Value = '  -1.57%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.799'
$ws.Range("E46").Value = '  +1.36%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.5946'
$ws.Range("E47").Value = '  -0.80%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.083'
$ws.Range("E48").Value = '  -1.62%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.279'
$ws.Range("E49").Value = '  +1.84%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '124.48'
$ws.Range("E50").Value = '  +0.54%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.07314'
$ws.Range("E51").Value = '  +0.31%  '
